# [ANV] updated decay chains spreadsheet
# Adds "Pb Density" and "Pb Target Fractions" worksheets (modeled on the
# existing "Cu Density" / "Cu Target Fractions " sheets), and nudges a
# handful of saved cell-selections on other sheets, matching the upstream
# commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "Pb Density" -- copy of "Cu Density", values updated for Pb
# ---------------------------------------------------------------------
$cuDensity = $wb.Worksheets.Item("Cu Density")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cuDensity.Copy($null, $lastSheet)
$pbDensity = $wb.Worksheets.Item($wb.Worksheets.Count)
$pbDensity.Name = "Pb Density"

$pbDensity.Range("A2").Value = "Pb"
$pbDensity.Range("B2").Value = 82
$pbDensity.Range("C2").Value = 1
$pbDensity.Range("D2").Value = 11340
$pbDensity.Range("G2").Value = 207.2
$pbDensity.Range("A6").Value = "https://en.wikipedia.org/wiki/Lead"
$pbDensity.Range("D15").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) New sheet "Pb Target Fractions" -- copy of "Cu Target Fractions ",
#    extended with Pb's four stable isotopes (204/206/207/208)
# ---------------------------------------------------------------------
$cuTargetFractions = $wb.Worksheets.Item("Cu Target Fractions ")
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$cuTargetFractions.Copy($null, $lastSheet2)
$pbTargetFractions = $wb.Worksheets.Item($wb.Worksheets.Count)
$pbTargetFractions.Name = "Pb Target Fractions"

# Cu's sheet only has 2 isotope rows (3:4); Pb needs 4 (3:6) -- insert two
# extra rows above the existing "Total" row so it lands on row 7.
$pbTargetFractions.Range("A5:A6").EntireRow.Insert()

$pbTargetFractions.Range("A2").Value = 82
$pbTargetFractions.Range("B2").Value = 1

$pbTargetFractions.Range("B3").Value = 1
$pbTargetFractions.Range("C3").Value = 204
$pbTargetFractions.Range("D3").Formula = "=1.4*(1-0.06)"
$pbTargetFractions.Range("E3").Formula = "=1.4*(1+0.06)"

$pbTargetFractions.Range("B4").Value = 1
$pbTargetFractions.Range("C4").Value = 206
$pbTargetFractions.Range("D4").Formula = "=24.1*(0.7)"
$pbTargetFractions.Range("E4").Formula = "=24.1*1.3"

$pbTargetFractions.Range("B5").Value = 1
$pbTargetFractions.Range("C5").Value = 207
$pbTargetFractions.Range("D5").Formula = "=22.1*0.5"
$pbTargetFractions.Range("E5").Formula = "=22.1*1.5"

$pbTargetFractions.Range("B6").Value = 1
$pbTargetFractions.Range("C6").Value = 208
$pbTargetFractions.Range("D6").Formula = "=52.4-0.7"
$pbTargetFractions.Range("E6").Formula = "=52.4+0.7"

$pbTargetFractions.Range("F3:F6").Formula = "=B3*(D3+E3)/200"
$pbTargetFractions.Range("G3:G6").Formula = "=B3*D3/100"
$pbTargetFractions.Range("H3:H6").Formula = "=B3*E3/100"

$pbTargetFractions.Range("H14").Select()

# This is the sheet that ends up with focus/tabSelected.
$pbTargetFractions.Activate()

# ---------------------------------------------------------------------
# 3) Saved-selection nudges on pre-existing sheets (unrelated to Pb, but
#    part of the same upstream commit).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("HDPE Target Fractions").Range("F4").Select()
$wb.Worksheets.Item("Shotcrete Target Fractions").Range("F31").Select()
$wb.Worksheets.Item("HDPE Density").Range("F31").Select()

# Leave the user back on the newly added / tab-selected sheet.
$pbTargetFractions.Activate()
